# Generate Report for Handback
# Updates the localization-status report after a handback was received for
# the "26dc01cb-1462-491f-a680-e2068927e702.md" file: its status flips from
# "In Translation" to "Handed back: in sync with en-US", its Latest Handback
# DateTime is refreshed, and the stale "handback file is not the latest"
# error is cleared.

$wb = $excel.ActiveWorkbook

$oldStatus = "In Translation"
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: zh-cn / de-de status columns for the 26dc01cb... row ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# --- zh-cn sheet: row 2 is the 26dc01cb-...md file ---
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("L2").Value = "2017-02-22 08:20:14"
$wsZhCn.Range("R2").Value = ""

# --- de-de sheet: row 2 is the 26dc01cb-...md file ---
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("L2").Value = "2017-02-22 08:20:35"
$wsDeDe.Range("R2").Value = ""

# --- Resize the Status columns that now hold the longer text ---
$statusColWidth = 29.17
$wsOverview.Columns.Item(5).ColumnWidth = $statusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $statusColWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $statusColWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $statusColWidth
